$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stahl_Treppe")
$ws.Rows.Item(9).Delete()
$ws.Select()
$ws.Range("L18").Select()
